$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '60.966.69'
$ws.Cells.Item(2, 5).Value = '  +2.84%  '
$ws.Cells.Item(3, 4).Value = '2.613.24'
$ws.Cells.Item(3, 5).Value = '  +1.52%  '
$ws.Cells.Item(4, 5).Value = '  +0.00%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '582.78'
$ws.Cells.Item(5, 5).Value = '  +4.76%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '144.53'
$ws.Cells.Item(6, 5).Value = '  +2.11%  '
$ws.Cells.Item(8, 5).Value = '  +0.95%  '
$ws.Cells.Item(9, 4).Value = '2.639.46'
$ws.Cells.Item(9, 5).Value = '  +2.34%  '
$ws.Cells.Item(10, 5).Value = '  -3.35%  '
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.107'
$ws.Cells.Item(11, 5).Value = '  +2.59%  '
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '0.159'
$ws.Cells.Item(12, 5).Value = '  -3.79%  '
$ws.Cells.Item(13, 5).Value = '  +6.48%  '
$ws.Cells.Item(14, 4).Value = '3.078.17'
$ws.Cells.Item(14, 5).Value = '  +1.70%  '
$ws.Cells.Item(15, 4).Value = '60.937.22'
$ws.Cells.Item(15, 5).Value = '  +2.82%  '
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '23.51'
$ws.Cells.Item(16, 5).Value = '  +2.61%  '
$ws.Cells.Item(17, 5).Value = '  +4.42%  '
$ws.Cells.Item(18, 4).Value = '2.627.24'
$ws.Cells.Item(18, 5).Value = '  +2.07%  '
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '11.32'
$ws.Cells.Item(19, 5).Value = '  +9.74%  '
$ws.Cells.Item(20, 5).Value = '  +3.17%  '
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '350.96'
$ws.Cells.Item(21, 5).Value = '  +3.90%  '
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '6.99'
$ws.Cells.Item(22, 5).Value = '  +8.20%  '
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '0.999'
$ws.Cells.Item(23, 5).Value = '  +0.13%  '
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '0.518'
$ws.Cells.Item(24, 5).Value = '  +7.83%  '
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '63.39'
$ws.Cells.Item(25, 5).Value = '  +1.16%  '
$ws.Cells.Item(26, 5).Value = '  +0.00%  '
$ws.Cells.Item(27, 5).Value = '  +0.65%  '
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '7.95'
$ws.Cells.Item(28, 5).Value = '  +7.59%  '
$ws.Cells.Item(29, 4).Value = '0.0₃0804'
$ws.Cells.Item(29, 5).Value = '  +4.08%  '
$ws.Cells.Item(30, 5).Value = '  +9.36%  '
$ws.Cells.Item(31, 5).Value = '  -0.14%  '
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '6.36'
$ws.Cells.Item(32, 5).Value = '  +2.85%  '
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '162.90'
$ws.Cells.Item(33, 5).Value = '  +2.35%  '
$ws.Cells.Item(34, 5).Value = '  +2.81%  '
$ws.Cells.Item(35, 2).Value = 'Fetch.AI'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '1.02'
$ws.Cells.Item(35, 5).Value = '  +13.68%  '
$ws.Cells.Item(36, 2).Value = 'NEARProtocol'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '4.32'
$ws.Cells.Item(36, 5).Value = '  +5.93%  '
$ws.Cells.Item(37, 5).Value = '  +6.72%  '
$ws.Cells.Item(38, 5).Value = '  +10.59%  '
$ws.Cells.Item(39, 2).Value = 'Filecoin'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '3.93'
$ws.Cells.Item(39, 5).Value = '  +7.07%  '
$ws.Cells.Item(40, 2).Value = 'OKB'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '37.99'
$ws.Cells.Item(40, 5).Value = '  +1.54%  '
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '310.28'
$ws.Cells.Item(41, 5).Value = '  +6.94%  '
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '0.850'
$ws.Cells.Item(42, 5).Value = '  -0.21%  '
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '134.58'
$ws.Cells.Item(43, 5).Value = '  -2.79%  '
$ws.Cells.Item(44, 2).Value = 'EnergySwap'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '19.98'
$ws.Cells.Item(44, 5).Value = '  +6.02%  '
$ws.Cells.Item(45, 2).Value = 'RenderToken'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '5.06'
$ws.Cells.Item(45, 5).Value = '  +11.93%  '
$ws.Cells.Item(46, 2).Value = 'Mantle'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '0.611'
$ws.Cells.Item(46, 5).Value = '  +3.15%  '
$ws.Cells.Item(47, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '20.43'
$ws.Cells.Item(47, 5).Value = '  +9.52%  '
$ws.Cells.Item(48, 2).Value = 'FirstDigitalUSD'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '0.995'
$ws.Cells.Item(48, 5).Value = '  -0.41%  '
$ws.Cells.Item(49, 2).Value = 'Stellar'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '0.0985'
$ws.Cells.Item(49, 5).Value = '  +1.41%  '
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '0.0553'
$ws.Cells.Item(50, 5).Value = '  +4.29%  '
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '0.0244'
$ws.Cells.Item(51, 5).Value = '  +4.41%  '
